$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.901.60"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "3.169.86"
$ws.Range("E3").Value = "  -4.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'590.79"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").Value = "'134.51"
$ws.Range("E6").Value = "  -4.97%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.167.07"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  -6.52%  "
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = "  -6.25%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  -4.90%  "
$ws.Range("D14").Value = "'34.28"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "3.690.63"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "3.174.72"
$ws.Range("E17").Value = "  -4.03%  "
$ws.Range("D18").Value = "62.876.17"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "'6.52"
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").Value = "'459.08"
$ws.Range("E20").Value = "  -4.59%  "
$ws.Range("D21").Value = "'13.95"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'0.695"
$ws.Range("E22").Value = "  -6.26%  "
$ws.Range("D23").Value = "'7.57"
$ws.Range("E23").Value = "  -5.41%  "
$ws.Range("D24").Value = "'13.25"
$ws.Range("E24").Value = "  -4.85%  "
$ws.Range("D25").Value = "'82.25"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'2.66"
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("D29").Value = "'6.73"
$ws.Range("D30").Value = "'7.61"
$ws.Range("E30").Value = "  -7.08%  "
$ws.Range("D31").Value = "'2.02"
$ws.Range("E31").Value = "  -5.87%  "
$ws.Range("D32").Value = "'27.09"
$ws.Range("E32").Value = "  -6.13%  "
$ws.Range("D33").Value = "'0.101"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "'2.35"
$ws.Range("E34").Value = "  -6.78%  "
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("D36").Value = "'5.78"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("D37").Value = "'51.21"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "0.0₃0704"
$ws.Range("E38").Value = "  -6.04%  "
$ws.Range("D39").Value = "'0.0386"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "'402.13"
$ws.Range("E40").Value = "  -7.36%  "
$ws.Range("D41").Value = "'8.07"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("D43").Value = "'0.111"
$ws.Range("E43").Value = "  -6.58%  "
$ws.Range("D44").Value = "2.812.63"
$ws.Range("E44").Value = "  -9.67%  "
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("D48").Value = "'124.09"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'34.62"
$ws.Range("E49").Value = "  -6.19%  "
$ws.Range("D50").Value = "'25.02"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("E51").Value = "  -2.49%  "
